$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17-28 shift down to 18-29.
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new weekly price record.
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 45216
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = "Chirimoya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 21000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2100
$ws.Range("T17").Value = 10
